$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Agregue mi correo": append Emmanuel's e-mail address to the end
#    of his line in the "Integrantes" list (2nd paragraph).
# ------------------------------------------------------------------
$emmanuel = $d.Paragraphs(2)
$emmanuel.Range.InsertAfter(" – rootandtoor@hotmail.com")

# ------------------------------------------------------------------
# 2) "borre documento temporal de word": remove the stray _GoBack
#    bookmark that Word leaves behind to mark the last edit point -
#    a throwaway / "temporary" artifact that shouldn't ship in the
#    saved document.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Re-saving the file with a current Word build also pulls in the
#    (unused) built-in "Heading 2" style definition/gallery entry.
#    Touch the style once so it gets minted into styles.xml, matching
#    that side effect, then leave the paragraph itself as "Normal".
# ------------------------------------------------------------------
$tmpPara = $d.Paragraphs(2)
$tmpPara.Style = "Heading 2"
$tmpPara.Style = "Normal"

Write-Output "done"
